$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 40.018403858124941
$ws.Range("C2").Value = 23.269420828125021
$ws.Range("D2").Value = 34.490713858124934
$ws.Range("E2").Value = 31.556935828125006

# Row 3 values
$ws.Range("B3").Value = 34.57344337125005
$ws.Range("C3").Value = 20.167593688124953
$ws.Range("D3").Value = 29.026893371250083
$ws.Range("E3").Value = 22.914452188124926

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
